$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "SELECT DISTINCT `n    prg.program_name AS `"Program`",`n  CASE`n    WHEN prg.program_link IS NOT NULL THEN prg.program_acronym`n        ELSE prg.program_link`n    END  AS `"Website`",`n    prg.focus_area AS `"Special Topic`",`n    prg.cancer_type AS `"Cancer Type`",`n CASE `n        WHEN prg.data_link IS NOT NULL THEN prg.program_acronym    `n        ELSE prg.data_link`n    END AS `"Data Location Details`"`nFROM `n    df_program prg`nWHERE `n    prg.focus_area LIKE '%Health Disparities%'`nORDER BY `n    lower(prg.program_name) ASC`nLIMIT 100;"
$ws.Range("B3").Value = "SELECT DISTINCT`n    prj.project_id AS `"Project ID`", `n    prj.project_title AS `"Project Title`",`n    prj.project_org_name AS `"Organization`",`n    prj.project_start_date AS `"Project Start Date`",`n    prj.project_end_date AS `"Project End Date`"`nFROM `n    df_project prj`nLEFT JOIN `n    df_program prg ON prj.`"program.program_id`" = prg.program_id`nLEFT JOIN `n    df_grant gnt ON prj.project_id = gnt.`"project.project_id`"`nLEFT JOIN `n    df_publication pub ON prj.project_id = pub.`"project.project_id`"`nWHERE `n    prg.focus_area LIKE '%Health Disparities%'`nORDER BY `n    lower(prj.project_id) ASC`nLIMIT 100;"
$ws.Range("B4").Value = "SELECT DISTINCT`n    gnt.grant_id AS `"Grant ID`", `n    prj.project_id AS `"Project`",`n    gnt.grant_title AS `"Grant Title`",`n    gnt.principal_investigators AS `"Principal Investigators`",`n    gnt.program_officers AS `"Program Officers`",`n    gnt.fiscal_year AS `"Fiscal Year`",`n    gnt.grant_end_date AS `"Project End Date`"`nFROM `n    df_grant gnt`nLEFT JOIN `n    df_project prj ON gnt.`"project.project_id`" = prj.project_id`nLEFT JOIN `n    df_program prg ON prj.`"program.program_id`" = prg.program_id`nLEFT JOIN `n    df_publication pub ON prj.project_id = pub.`"project.project_id`"`nWHERE `n    prg.focus_area LIKE '%Health Disparities%'`nORDER BY `n    lower(gnt.grant_id) ASC`nLIMIT 100;"
$ws.Range("B5").Value = "SELECT DISTINCT`n    pub.pmid AS `"PubMed ID`", `n    pub.publication_title AS `"Title`",`n    pub.authors AS `"Authors`",`n    pub.publication_date AS `"Publication Date`",`n    pub.cited_by AS `"Cited By`",`n    CASE `n    WHEN pub.relative_citation_ratio = 0 THEN '0'`n    WHEN pub.relative_citation_ratio = 7.0 THEN '7'`nWHEN pub.relative_citation_ratio = 1.0 THEN '1'`n    WHEN pub.relative_citation_ratio = ROUND(pub.relative_citation_ratio) THEN CAST(ROUND(pub.relative_citation_ratio) AS VARCHAR) `n    ELSE CAST(ROUND(pub.relative_citation_ratio, 2) AS VARCHAR)`nEND AS `"Relative Citation Ratio`"`nFROM `n    df_publication pub`nLEFT JOIN `n    df_project prj ON pub.`"project.project_id`" = prj.project_id`nLEFT JOIN `n    df_program prg ON prj.`"program.program_id`" = prg.program_id`nLEFT JOIN `n    df_grant gnt ON prj.project_id = gnt.`"project.project_id`"`nWHERE `n    prg.focus_area LIKE '%Health Disparities%'`nORDER BY `n    lower(pub.pmid) ASC`nLIMIT 100;"

$ws.Range("B5").Select()
